# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
# Bulgaria First League.xlsx update

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rows 11 and 12 describe two Levski Sofia / CSKA 1948 Sofia fixtures
#    played on the same date whose results/odds had been mixed up.
#    Correct them: row 11 becomes the former row-12 fixture data (and
#    vice-versa), while the id (A) and date (E) stay anchored to their
#    original row.
# ---------------------------------------------------------------------

# --- Row 11 final values (previously held by row 12) ---
$ws.Range("B11").Value = 6627724
$ws.Range("F11").Value = "CSKA 1948 Sofia"
$ws.Range("G11").Value = "Lokomotiv Plovdiv"
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "D"
$ws.Range("K11").Value = 1.5
$ws.Range("L11").Value = 3.8
$ws.Range("M11").Value = 6
$ws.Range("N11").Value = 1.45
$ws.Range("O11").Value = 4.2
$ws.Range("P11").Value = 8
$ws.Range("Q11").Value = -1.25
$ws.Range("R11").Value = 2.025
$ws.Range("S11").Value = 1.825
$ws.Range("T11").Value = 2.5
$ws.Range("U11").Value = 1.85
$ws.Range("V11").Value = 2
$ws.Range("W11").Value = -1
$ws.Range("X11").Value = 3.2
$ws.Range("Y11").Value = -1
$ws.Range("Z11").Value = -1
$ws.Range("AA11").Value = 0.825
$ws.Range("AB11").Value = -1
$ws.Range("AC11").Value = 1

# --- Row 12 final values (previously held by row 11) ---
$ws.Range("B12").Value = 6627725
$ws.Range("F12").Value = "Levski Sofia"
$ws.Range("G12").Value = "CSKA Sofia"
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = "A"
$ws.Range("K12").Value = 2.625
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 2.6
$ws.Range("N12").Value = 2.55
$ws.Range("O12").Value = 3.3
$ws.Range("P12").Value = 2.8
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 1.825
$ws.Range("S12").Value = 2.025
$ws.Range("T12").Value = 2.25
$ws.Range("U12").Value = 1.85
$ws.Range("V12").Value = 2
$ws.Range("W12").Value = -1
$ws.Range("X12").Value = -1
$ws.Range("Y12").Value = 1.8
$ws.Range("Z12").Value = -1
$ws.Range("AA12").Value = 1.025
$ws.Range("AB12").Value = -0.5
$ws.Range("AC12").Value = 0.5

# ---------------------------------------------------------------------
# 2) Row 188 (Beroe vs Ludogorets Razgrad) was not played yet before;
#    fill in the final score / result plus refreshed closing odds.
# ---------------------------------------------------------------------
$ws.Range("H188").Value = 0
$ws.Range("I188").Value = 2
$ws.Range("J188").Value = "A"
$ws.Range("O188").Value = 5.75
$ws.Range("P188").Value = 1.222
$ws.Range("Q188").Value = 1.75
$ws.Range("R188").Value = 1.925
$ws.Range("S188").Value = 1.925
$ws.Range("U188").Value = 1.825
$ws.Range("V188").Value = 2.025
$ws.Range("W188").Value = -1
$ws.Range("X188").Value = -1
$ws.Range("Y188").Value = 0.222
$ws.Range("Z188").Value = -0.5
$ws.Range("AA188").Value = 0.4625
$ws.Range("AB188").Value = -1
$ws.Range("AC188").Value = 1.025

# ---------------------------------------------------------------------
# 3) Append a new fixture row 193 (CSKA Sofia vs Beroe), not yet played.
#    Copy formatting from the previous data row (192) first so the id
#    and date cells keep their usual number formats/styles. Only the
#    cells that actually need special formatting (A: bold/bordered id
#    style, E: date style) are copied, so we don't create stray empty
#    cells for the not-yet-applicable H/I/J score/result columns.
# ---------------------------------------------------------------------
$ws.Range("A192").Copy()
$ws.Range("A193").PasteSpecial(-4122)
$ws.Range("E192").Copy()
$ws.Range("E193").PasteSpecial(-4122)

$ws.Range("A193").Value = 191
$ws.Range("B193").Value = 6978410
$ws.Range("C193").Value = "Bulgaria First League"
$ws.Range("D193").Value = "Bulgaria First League"
$ws.Range("E193").Value = 45354.39583333334
$ws.Range("F193").Value = "CSKA Sofia"
$ws.Range("G193").Value = "Beroe"
$ws.Range("K193").Value = 1.222
$ws.Range("L193").Value = 6
$ws.Range("M193").Value = 12
$ws.Range("N193").Value = 1.25
$ws.Range("O193").Value = 5.75
$ws.Range("P193").Value = 9.5
$ws.Range("Q193").Value = -1.5
$ws.Range("R193").Value = 1.875
$ws.Range("S193").Value = 1.975
$ws.Range("T193").Value = 2.5
$ws.Range("U193").Value = 1.975
$ws.Range("V193").Value = 1.875
$ws.Range("W193").Value = 0
$ws.Range("X193").Value = 0
$ws.Range("Y193").Value = 0
$ws.Range("Z193").Value = 0
$ws.Range("AA193").Value = 0
